$wb = $excel.ActiveWorkbook

# Overview sheet: zh-cn / de-de status columns (E2:F3) go from "Ready for handoff" to "In Translation"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: Status column (C) for data rows
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: Status column (C) for data rows
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# The shorter status text means the status columns no longer need to be as
# wide - narrow the affected columns on all three sheets to match.
$wsOverview.Columns.Item(5).ColumnWidth = 13.4101845877511
$wsOverview.Columns.Item(6).ColumnWidth = 13.4101845877511
$wsZhCn.Columns.Item(3).ColumnWidth = 13.4101845877511
$wsDeDe.Columns.Item(3).ColumnWidth = 13.4101845877511
